$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits (shared-string content changes) ---
# "TYPO3 Wordpress" -> "Wordpress TYPO3" (only used by B35)
$ws.Range("B35").Value = "Wordpress TYPO3"

# "IPMA ICB4 Level D" / "Practice, People, Perspective" split differently:
#   -> "IPMA ICB4" / "Level D, Practice, People, Perspective"
# These shared strings are used across B36:B42 and C36:C42 - update every cell
# that referenced them so every usage moves together (same as editing the
# shared string table directly).
$rows = 36,37,38,39,40,41,42
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "IPMA ICB4"
    $ws.Range("C$r").Value = "Level D, Practice, People, Perspective"
}

# --- Column layout: split the combined A:B column into two distinct widths ---
# Column A keeps its existing (untouched) width; only column B is widened, which
# makes the sheet emit separate <col> entries for A and B.
$ws.Columns(2).ColumnWidth = 22

# --- Row height adjustments (shrinking rows whose wrapped text now fits tighter) ---
$ws.Rows(8).RowHeight = 46.7
$ws.Rows(11).RowHeight = 42.6
$ws.Rows(35).RowHeight = 56.45
$ws.Rows(36).RowHeight = 42.6
